# Swap the contents of columns B:AD between the given pairs of rows.
# Column A (the running index) is left untouched on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 2")

# Row pairs whose B:AD contents need to be swapped with each other.
$rowPairs = @(
    @(24, 25),
    @(66, 67),
    @(78, 79),
    @(84, 85),
    @(118, 119),
    @(160, 161),
    @(193, 194)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # NOTE: build the A1 addresses with string concatenation rather than
    # "B$r1:AD$r1" interpolation -- a bare ":" right after an interpolated
    # variable name is parsed as a scope/drive qualifier (like $env:PATH),
    # which silently turns "AD$r1" into an empty lookup.
    $addr1 = "B" + $r1 + ":AD" + $r1
    $addr2 = "B" + $r2 + ":AD" + $r2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
